$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.39185424165248
$ws.Range("C2").Value = 7.686111888090331
$ws.Range("D2").Value = 10.72449922999026
$ws.Range("F2").Value = 30.62284562759701
$ws.Range("G2").Value = 3.631325844477329
$ws.Range("I2").Value = 20.07958466336632
$ws.Range("J2").Value = 11.07894246573618
$ws.Range("M2").Value = 17.23771168851673
$ws.Range("N2").Value = 17.34590629550259
$ws.Range("O2").Value = 22.18463177278219
$ws.Range("B3").Value = 11.89855794563825
$ws.Range("C3").Value = 7.254239078304299
$ws.Range("D3").Value = 10.71056548528108
$ws.Range("F3").Value = 30.63581898255675
$ws.Range("G3").Value = 3.633359456668868
$ws.Range("I3").Value = 20.17228674823479
$ws.Range("J3").Value = 11.10574437791838
$ws.Range("M3").Value = 17.08556223724025
$ws.Range("N3").Value = 17.39539446391844
$ws.Range("O3").Value = 22.22826816596688
$ws.Range("B4").Value = 11.58628021379759
$ws.Range("C4").Value = 6.975710198148732
$ws.Range("D4").Value = 10.70366493595778
$ws.Range("F4").Value = 30.65203871426676
$ws.Range("G4").Value = 3.634675013304309
$ws.Range("I4").Value = 20.23364171717723
$ws.Range("J4").Value = 11.12376066836671
$ws.Range("M4").Value = 16.99425626692976
$ws.Range("N4").Value = 17.42753438375512
$ws.Range("O4").Value = 22.260679529988
$ws.Range("B5").Value = 11.45685013710231
$ws.Range("C5").Value = 6.858940037595769
$ws.Range("D5").Value = 10.7012713753404
$ws.Range("F5").Value = 30.66072179557056
$ws.Range("G5").Value = 3.635227991699058
$ws.Range("I5").Value = 20.25975771534744
$ws.Range("J5").Value = 11.13149467258511
$ws.Range("M5").Value = 16.9576144698933
$ws.Range("N5").Value = 17.4410736296537
$ws.Range("O5").Value = 22.27529648191869
$ws.Range("B6").Value = 11.43523294765226
$ws.Range("C6").Value = 6.839356151149289
$ws.Range("D6").Value = 10.70089926477046
$ws.Range("F6").Value = 30.66228875677062
$ws.Range("G6").Value = 3.635320834317437
$ws.Range("I6").Value = 20.26416144272277
$ws.Range("J6").Value = 11.13280258680993
$ws.Range("M6").Value = 16.95156530678636
$ws.Range("N6").Value = 17.44334853202478
$ws.Range("O6").Value = 22.27780860125097
$ws.Range("B7").Value = 11.58454321494165
$ws.Range("C7").Value = 6.974148489609481
$ws.Range("D7").Value = 10.70363095827207
$ws.Range("F7").Value = 30.65214742595188
$ws.Range("G7").Value = 3.634682402553777
$ws.Range("I7").Value = 20.23398942138519
$ws.Range("J7").Value = 11.1238633837875
$ws.Range("M7").Value = 16.99375976586759
$ws.Range("N7").Value = 17.42771518803348
$ws.Range("O7").Value = 22.2608709593887
$ws.Range("B8").Value = 12.22382064258246
$ws.Range("C8").Value = 7.54002360998839
$ws.Range("D8").Value = 10.71935289558261
$ws.Range("F8").Value = 30.62560479282609
$ws.Range("G8").Value = 3.632013178745095
$ws.Range("I8").Value = 20.11062654811917
$ws.Range("J8").Value = 11.08786002606831
$ws.Range("M8").Value = 17.18483266697161
$ws.Range("N8").Value = 17.362606324505
$ws.Range("O8").Value = 22.19850936964129
$ws.Range("B9").Value = 13.39576486480181
$ws.Range("C9").Value = 8.540573650021331
$ws.Range("D9").Value = 10.76320132135263
$ws.Range("F9").Value = 30.63909025749027
$ws.Range("G9").Value = 3.627307318243621
$ws.Range("I9").Value = 19.90399574661406
$ws.Range("J9").Value = 11.02963414415571
$ws.Range("M9").Value = 17.57474788715921
$ws.Range("N9").Value = 17.24880385492018
$ws.Range("O9").Value = 22.12094791641611
$ws.Range("B10").Value = 14.1987992094134
$ws.Range("C10").Value = 9.205924061828483
$ws.Range("D10").Value = 10.80318833486972
$ws.Range("F10").Value = 30.68893289139587
$ws.Range("G10").Value = 3.624168703962647
$ws.Range("I10").Value = 19.7738170615515
$ws.Range("J10").Value = 10.99440081529994
$ws.Range("M10").Value = 17.868369367455
$ws.Range("N10").Value = 17.17359513058876
$ws.Range("O10").Value = 22.09140944890998
$ws.Range("B11").Value = 14.5501515765068
$ws.Range("C11").Value = 9.493006063921792
$ws.Range("D11").Value = 10.82302695562068
$ws.Range("F11").Value = 30.72024772556085
$ws.Range("G11").Value = 3.622809365775435
$ws.Range("I11").Value = 19.71931987698845
$ws.Range("J11").Value = 10.98001048362731
$ws.Range("M11").Value = 18.00303685243608
$ws.Range("N11").Value = 17.14119330226017
$ws.Range("O11").Value = 22.08395900547281
$ws.Range("B12").Value = 14.68109931072249
$ws.Range("C12").Value = 9.599447566299235
$ws.Range("D12").Value = 10.83077254266322
$ws.Range("F12").Value = 30.73334341789369
$ws.Range("G12").Value = 3.622304405662317
$ws.Range("I12").Value = 19.69936456297205
$ws.Range("J12").Value = 10.97479668571287
$ws.Range("M12").Value = 18.05415083883964
$ws.Range("N12").Value = 17.12918314093406
$ws.Range("O12").Value = 22.0819999461647
$ws.Range("B13").Value = 14.65299222696838
$ws.Range("C13").Value = 9.576624835735329
$ws.Range("D13").Value = 10.82909408783668
$ws.Range("F13").Value = 30.73046808738226
$ws.Range("G13").Value = 3.622412723140592
$ws.Range("I13").Value = 19.70363194139174
$ws.Range("J13").Value = 10.9759090956142
$ws.Range("M13").Value = 18.04313789147123
$ws.Range("N13").Value = 17.13175820458282
$ws.Range("O13").Value = 22.08238349946012
$ws.Range("B14").Value = 14.56096734689405
$ws.Range("C14").Value = 9.501808708360853
$ws.Range("D14").Value = 10.82365954056822
$ws.Range("F14").Value = 30.72130035027147
$ws.Range("G14").Value = 3.622767626440356
$ws.Range("I14").Value = 19.7176644662704
$ws.Range("J14").Value = 10.9795768212234
$ws.Range("M14").Value = 18.00723988493356
$ws.Range("N14").Value = 17.14020001812498
$ws.Range("O14").Value = 22.08378054561623
$ws.Range("B15").Value = 14.5043231356106
$ws.Range("C15").Value = 9.455685289902057
$ws.Range("D15").Value = 10.82036096446129
$ws.Range("F15").Value = 30.71584582174886
$ws.Range("G15").Value = 3.622986288691813
$ws.Range("I15").Value = 19.72634863340659
$ws.Range("J15").Value = 10.98185408012197
$ws.Range("M15").Value = 17.9852655602998
$ws.Range("N15").Value = 17.14540467355898
$ws.Range("O15").Value = 22.08474859979499
$ws.Range("B16").Value = 14.17554827591897
$ws.Range("C16").Value = 9.186846287189624
$ws.Range("D16").Value = 10.80192467532832
$ws.Range("F16").Value = 30.68705986789979
$ws.Range("G16").Value = 3.624258912776105
$ws.Range("I16").Value = 19.7774738165643
$ws.Range("J16").Value = 10.99537421148709
$ws.Range("M16").Value = 17.85958723726456
$ws.Range("N16").Value = 17.1757490468618
$ws.Range("O16").Value = 22.09201695204336
$ws.Range("B17").Value = 13.97020970757411
$ws.Range("C17").Value = 9.01790695990929
$ws.Range("D17").Value = 10.79103408929396
$ws.Range("F17").Value = 30.67161061235363
$ws.Range("G17").Value = 3.625057118928847
$ws.Range("I17").Value = 19.81004866007963
$ws.Range("J17").Value = 11.00408780691935
$ws.Range("M17").Value = 17.78274063621599
$ws.Range("N17").Value = 17.19482766386559
$ws.Range("O17").Value = 22.09801035832514
$ws.Range("B18").Value = 13.85079593913478
$ws.Range("C18").Value = 8.919271408724507
$ws.Range("D18").Value = 10.78492550906788
$ws.Range("F18").Value = 30.66353835387682
$ws.Range("G18").Value = 3.625522670147943
$ws.Range("I18").Value = 19.82922912473204
$ws.Range("J18").Value = 11.00925376307626
$ws.Range("M18").Value = 17.73864617623861
$ws.Range("N18").Value = 17.20597168053744
$ws.Range("O18").Value = 22.10202106182024
$ws.Range("B19").Value = 13.81014289898076
$ws.Range("C19").Value = 8.885624281964171
$ws.Range("D19").Value = 10.78288405311615
$ws.Range("F19").Value = 30.66094512116352
$ws.Range("G19").Value = 3.625681406076198
$ws.Range("I19").Value = 19.83579952642934
$ws.Range("J19").Value = 11.01102933698604
$ws.Range("M19").Value = 17.72373590478331
$ws.Range("N19").Value = 17.20977415882178
$ws.Range("O19").Value = 22.10347573993516
$ws.Range("B20").Value = 13.99220450438639
$ws.Range("C20").Value = 9.036042805671302
$ws.Range("D20").Value = 10.7921773552519
$ws.Range("F20").Value = 30.67317102525972
$ws.Range("G20").Value = 3.624971481931425
$ws.Range("I20").Value = 19.80653501264407
$ws.Range("J20").Value = 11.00314427795306
$ws.Range("M20").Value = 17.79091043798241
$ws.Range("N20").Value = 17.19277907153802
$ws.Range("O20").Value = 22.09731402277676
$ws.Range("B21").Value = 14.5880550202592
$ws.Range("C21").Value = 9.523845825500279
$ws.Range("D21").Value = 10.8252495039171
$ws.Range("F21").Value = 30.72395960095061
$ws.Range("G21").Value = 3.622663117379957
$ws.Range("I21").Value = 19.71352425460513
$ws.Range("J21").Value = 10.97849312943495
$ws.Range("M21").Value = 18.0177811020556
$ws.Range("N21").Value = 17.1377134100195
$ws.Range("O21").Value = 22.08334678961247
$ws.Range("B22").Value = 14.96519198449345
$ws.Range("C22").Value = 9.829411996220932
$ws.Range("D22").Value = 10.8482209645432
$ws.Range("F22").Value = 30.76436227025403
$ws.Range("G22").Value = 3.621211519161816
$ws.Range("I22").Value = 19.65671038343471
$ws.Range("J22").Value = 10.96375490924014
$ws.Range("M22").Value = 18.16672576822759
$ws.Range("N22").Value = 17.10323827845687
$ws.Range("O22").Value = 22.07924455246704
$ws.Range("B23").Value = 14.76505913967604
$ws.Range("C23").Value = 9.667544890203104
$ws.Range("D23").Value = 10.83583785471746
$ws.Range("F23").Value = 30.74214100586587
$ws.Range("G23").Value = 3.621981060277748
$ws.Range("I23").Value = 19.68666852447528
$ws.Range("J23").Value = 10.9714953607084
$ws.Range("M23").Value = 18.08718274363452
$ws.Range("N23").Value = 17.12150005256511
$ws.Range("O23").Value = 22.08097379302091
$ws.Range("B24").Value = 13.98226489132836
$ws.Range("C24").Value = 9.027848291165181
$ws.Range("D24").Value = 10.79166000920922
$ws.Range("F24").Value = 30.6724630399532
$ws.Range("G24").Value = 3.625010177693784
$ws.Range("I24").Value = 19.80812212252249
$ws.Range("J24").Value = 11.00357036024435
$ws.Range("M24").Value = 17.78721660071317
$ws.Range("N24").Value = 17.1937046934969
$ws.Range("O24").Value = 22.09762707618034
$ws.Range("B25").Value = 13.08838545350363
$ws.Range("C25").Value = 8.281961324599243
$ws.Range("D25").Value = 10.74996149751622
$ws.Range("F25").Value = 30.62842253273281
$ws.Range("G25").Value = 3.628524153107349
$ws.Range("I25").Value = 19.95610490546211
$ws.Range("J25").Value = 11.04406072692036
$ws.Range("M25").Value = 17.46785368604156
$ws.Range("N25").Value = 17.27811092245188
$ws.Range("O25").Value = 22.13712142376376
